# Applies the "add MI healthy climate plan" edit:
#  - Add explanatory "GRA for ..." summary lines to the About sheet (rows 54-65)
#  - Zero out the Deficit Spending and Payroll Taxes weights for Carbon Tax
#    Revenue, leaving only the Household Taxes weight of 5.

$wb = $excel.ActiveWorkbook

# --- About sheet: append the new GRA weighting summary notes --------------
$about = $wb.Worksheets.Item("About")

$notes = @(
    "GRA for Carbon Tax Revenue[household taxes] = 5",
    "GRA for Fuel Tax Revenue[household taxes] = 5",
    "GRA for EV Subsidy[deficit spending] = 5",
    "GRA for Electricity Generation Subsidies[deficit spending] = 5",
    "GRA for Electricity Capacity Construction Subsidies[deficit spending] = 5",
    "GRA for Distributed Solar Subsidy[deficit spending] = 5",
    "GRA for Fuel Subsidies[deficit spending] = 5",
    "GRA for National Debt Interest[regular spending] = 5",
    "GRA for National Debt Interest[household taxes] = 5",
    "GRA for Remaining Government Cash Flow Changes[regular spending] = 5",
    "GRA for Remaining Government Cash Flow Changes[deficit spending] = 5",
    "GRA for Remaining Government Cash Flow Changes[household taxes] = 5"
)

$startRow = 54
for ($i = 0; $i -lt $notes.Length; $i++) {
    $about.Cells.Item($startRow + $i, 1).Value = $notes[$i]
}

# --- Set Values Here: Carbon Tax Revenue row -------------------------------
# Row 8 = carbon tax revenue weights: B=Regular, C=Deficit, D=Household, E=Payroll, F=Corporate
# GRA-carbontax!B2:B6 is a TRANSPOSE() array formula that mirrors this row,
# so updating it here automatically ripples the Deficit/Payroll weights to 0
# on the GRA-carbontax sheet as well.
$setValues = $wb.Worksheets.Item("Set Values Here")
$setValues.Range("C8").Value = 0
$setValues.Range("E8").Value = 0

$excel.Calculate()
